# Weekly update: insert a new price record as row 7 (most recent week),
# pushing the previously-existing rows 7-10 down to rows 8-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(7, 3).Value = "La Araucanía"
$ws.Cells.Item(7, 4).Value = 44519
$ws.Cells.Item(7, 5).Value = 9
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100104
$ws.Cells.Item(7, 8).Value = "Frutos de pepita"
$ws.Cells.Item(7, 9).Value = 100104004
$ws.Cells.Item(7, 10).Value = "Níspero"
$ws.Cells.Item(7, 11).Value = "Californiana(o)"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 30
$ws.Cells.Item(7, 14).Value = 28000
$ws.Cells.Item(7, 15).Value = 28000
$ws.Cells.Item(7, 16).Value = 28000
$ws.Cells.Item(7, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(7, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(7, 19).Value = 2800
$ws.Cells.Item(7, 20).Value = 10
